$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = 5
$ws.Range("C7").Value = "get flash dump"
$ws.Range("E7").Value = "target board address"
$ws.Range("B7").Value = "dump"
$ws.Range("D7").Value = 1

$ws.Range("E7").Select()

$excel.ActiveWindow.WindowState = -4140
